$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-template_type")

# New code/description strings for "Part 4" of the Registration Acknowledgement Template
$code = "reg-ack-template-part4"
$descrEng = "Registration Acknowledgement Template - Part 4"
$descrAra = "نموذج شكر التسجيل"
$descrFra = "accusé de réception"

$startRow = 122

# Row 122: English
$ws.Cells.Item($startRow, 1).Value = $code
$ws.Cells.Item($startRow, 2).Value = $descrEng
$ws.Cells.Item($startRow, 3).Value = "eng"
$ws.Cells.Item($startRow, 4).Value = $true
$ws.Cells.Item($startRow, 5).Value = "superadmin"
$ws.Cells.Item($startRow, 6).Value = "now()"

# Row 123: Arabic
$ws.Cells.Item($startRow + 1, 1).Value = $code
$ws.Cells.Item($startRow + 1, 2).Value = $descrAra
$ws.Cells.Item($startRow + 1, 3).Value = "ara"
$ws.Cells.Item($startRow + 1, 4).Value = $true
$ws.Cells.Item($startRow + 1, 5).Value = "superadmin"
$ws.Cells.Item($startRow + 1, 6).Value = "now()"

# Row 124: French
$ws.Cells.Item($startRow + 2, 1).Value = $code
$ws.Cells.Item($startRow + 2, 2).Value = $descrFra
$ws.Cells.Item($startRow + 2, 3).Value = "fra"
$ws.Cells.Item($startRow + 2, 4).Value = $true
$ws.Cells.Item($startRow + 2, 5).Value = "superadmin"
$ws.Cells.Item($startRow + 2, 6).Value = "now()"

# Update selection to mirror the post-edit state (select the row below the data down to the end)
$ws.Range("A125:XFD1048576").Select()
